$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date value from 45224 to 45233 for rows 2-7
$ws.Range("C2:C7").Value = 45233
